# Update workbook for carjacking-by-neighborhood-by-month: refresh "through" date
# from 2022-03-04 to 2022-03-05 and add newly-reported incidents (commit: "Add
# data for 2022-03-13").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new cutoff date.
$ws.Name = "Through 2022-03-05"

# Update the column header text (row 1, column B) to match the new cutoff date.
$ws.Range("B1").Value = "March 2022 (through March 05)"

# Update the March-2022 counts (column B) and a few other incident counts
# that changed in other month columns, per neighborhood.

# Austin (row 3): March 2022 count 1 -> 2
$ws.Range("B3").Value = 2

# North Lawndale (row 4): March 2017 (Q) count 1 -> 2
$ws.Range("Q4").Value = 2

# Garfield Park (row 5): March 2022 count 2 -> 3
$ws.Range("B5").Value = 3

# Kenwood (row 8): new March 2018 (N) count of 1
$ws.Range("N8").Value = 1

# Chicago Lawn (row 10): new March 2022 (B) count of 2
$ws.Range("B10").Value = 2

# West Loop (row 11): March 2020 (H) count 1 -> 2
$ws.Range("H11").Value = 2

# Englewood (row 12): new March 2022 (B) count of 2; new March 2016 (T) count of 1
$ws.Range("B12").Value = 2
$ws.Range("T12").Value = 1

# Humboldt Park (row 15): new March 2022 (B) count of 1
$ws.Range("B15").Value = 1

# Grand Crossing (row 26): new March 2022 (B) count of 1
$ws.Range("B26").Value = 1

# Calumet Heights (row 27): new March 2018 (N) count of 1
$ws.Range("N27").Value = 1

# Logan Square (row 28): new March 2019 (K) count of 1
$ws.Range("K28").Value = 1

# New City (row 32): March 2022 count 1 -> 2
$ws.Range("B32").Value = 2

# Morgan Park (row 40): new March 2022 (B) count of 1
$ws.Range("B40").Value = 1

# Little Village (row 49): new March 2020 (H) count of 1
$ws.Range("H49").Value = 1

# Albany Park (row 50): new March 2018 (N) count of 1
$ws.Range("N50").Value = 1

# Clearing (row 54): March 2021 (E) count 1 -> 2; new March 2015 (W) count of 1
$ws.Range("E54").Value = 2
$ws.Range("W54").Value = 1

# Fuller Park (row 60): new March 2016 (T) count of 1
$ws.Range("T60").Value = 1
